# GMM.xlsx report update ("Improved report and translated"):
#   - split the "Tied Full-Cov 8 GAU" test-dataset row into two rows: the
#     original model (now with a trailing space in its label) and a new
#     "(balanced)" variant with its own error rate;
#   - widen column B so the longer labels fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch row 2 so it is materialised in the sheet (minor formatting nudge).
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(2).RowHeight

# Widen column B to fit the new, longer row labels.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth + 6

# Row 19 label gains a trailing space (it now has a sibling row below it).
$ws.Cells.Item(19, 2).Value = "Tied Full-Cov 8 GAU "

# New row 20: the "(balanced)" variant - same label style as row 19, same
# Min DCF / Act DCF as row 19, but its own (worse) error rate.
$ws.Cells.Item(19, 2).Copy()
$ws.Cells.Item(20, 2).PasteSpecial(-4122)
$ws.Cells.Item(20, 2).Value = "Tied Full-Cov 8 GAU (balanced) "

$ws.Cells.Item(19, 4).Copy($ws.Cells.Item(20, 4))
$ws.Cells.Item(19, 5).Copy($ws.Cells.Item(20, 5))

# F20 needs the literal text "15.0" (not the number 15) with default
# formatting. Build it via a throwaway text formula, then paste just the
# resulting value so no stray number format / quote-prefix style is left
# behind on the cell.
$ws.Cells.Item(500, 500).Formula = '="15.0"'
$ws.Cells.Item(500, 500).Copy()
$ws.Cells.Item(20, 6).PasteSpecial(-4163)
$ws.Cells.Item(500, 500).ClearContents()

$excel.CutCopyMode = $false
